$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 5 corresponds to the
# e88d80a3-4555-40bc-a626-2f3ac0edf367...zh-cn.xlf entry.
# Column D = Correspond Handoff Datetime, Column G = Correspond Handback DateTime
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D5").Value = "2016-01-26 06:04:55"
$wsZh.Range("G5").Value = "2016-01-26 06:05:57"

# de-de sheet: row 5 corresponds to the
# e88d80a3-4555-40bc-a626-2f3ac0edf367...de-de.xlf entry.
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D5").Value = "2016-01-26 06:05:10"
$wsDe.Range("G5").Value = "2016-01-26 06:06:20"
